# Transform the vertical label/value list (A1:B11) into a horizontal
# header row (row 1) + single data row (row 2) spanning columns A:K.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the old A1:B11 content so nothing stale remains outside the
# new A1:K2 dimension.
$ws.Range("A1:B11").Clear()

# New header row (row 1).
$ws.Cells.Item(1, 1).Value = "MIGRATION DATE"
$ws.Cells.Item(1, 2).Value = "FINANCIAL INSTITUTION NAME"
$ws.Cells.Item(1, 3).Value = "ENTITY ID"
$ws.Cells.Item(1, 4).Value = "ADDRESS"
$ws.Cells.Item(1, 5).Value = "CITY"
$ws.Cells.Item(1, 6).Value = "STATE"
$ws.Cells.Item(1, 7).Value = "ZIP CODE"
$ws.Cells.Item(1, 8).Value = "PHONE #"
$ws.Cells.Item(1, 9).Value = "PROJECT COORDINATOR"
$ws.Cells.Item(1, 10).Value = "CERTIFICATION REQUIRED (Yes or No)"
$ws.Cells.Item(1, 11).Value = "CERTIFICATION COORDINATOR"

# New data row (row 2). Force text formatting first on cells whose
# literal values look numeric/date-like so Excel keeps them as text
# instead of auto-converting to a date serial or a number.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2025-10-16"
$ws.Cells.Item(2, 2).Value = "YYY"
$ws.Cells.Item(2, 3).Value = "123ABX007"
$ws.Cells.Item(2, 4).Value = "Karapakkam"
$ws.Cells.Item(2, 5).Value = "Chennai"
$ws.Cells.Item(2, 6).Value = "Tamil Nadu"
$ws.Cells.Item(2, 7).NumberFormat = "@"
$ws.Cells.Item(2, 7).Value = "600117"
$ws.Cells.Item(2, 8).NumberFormat = "@"
$ws.Cells.Item(2, 8).Value = "9911991100"
$ws.Cells.Item(2, 9).Value = "Sam"
$ws.Cells.Item(2, 10).Value = "Yes"

# K2 mirrors the original blank CERTIFICATION COORDINATOR value cell
# (B11 was empty) - keep it present in the sheet but with no content.
$ws.Cells.Item(2, 11).NumberFormat = "@"
